# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# Source workbook layout (before):
#   Sheet1 "总计"     - summary table (A:D), rows 2-7 = 2022-Q2 .. 2021-Q1
#   Sheet2 "2022-Q2"
#   Sheet3 "2022-Q1"
#   Sheet4 "2021-Q4"
#   Sheet5 "2021-Q3"
#   Sheet6 "2021-Q2"
#   Sheet7 "2021-Q1"
#
# Target workbook layout (after):
#   Sheet1 "总计"      - same summary table, now with an extra 2022-Q3 row
#                         inserted right after the header (rows 2-8)
#   Sheet2 "2022-Q3"   - brand-new fund-holder detail sheet (A1:H23)
#   Sheet3 "2022-Q2"   - (was Sheet2, unchanged content, just shifted right)
#   Sheet4 "2022-Q1"   - (was Sheet3)
#   Sheet5 "2021-Q4"   - (was Sheet4)
#   Sheet6 "2021-Q3"   - (was Sheet5)
#   Sheet7 "2021-Q2"   - (was Sheet6)
#   Sheet8 "2021-Q1"   - (was Sheet7)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Excel enum constants (kept local so the script has no external deps)
$xlContinuous   = 1
$xlCenter       = -4108
$xlTop          = -4160
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while forcing Excel to keep it as TEXT
# (many of the source columns hold numeric-looking strings such as "87.00"
# or zero-padded fund codes such as "004997" that must not be coerced into
# numbers, which would drop trailing zeros / leading zeros).
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    # Flip the style back to the workbook default so no stray "@"-formatted
    # style lingers on the cell (matches how the rest of the sheet looks).
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Helper: populate one data row (row $r) of the new "2022-Q3" sheet.
# Column layout: A index(n) | B code(txt) | C name(txt) | D size(txt) |
#                E position(txt) | F ratio(txt) | G value(txt) | H rank(n)
# ---------------------------------------------------------------------------
function Set-Q3Row($ws, $r, $a, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$r").Value2 = $a

    $txtRng = $ws.Range("B$r" + ":G$r")
    $txtRng.NumberFormat = "@"
    $ws.Range("B$r").Value2 = $b
    $ws.Range("C$r").Value2 = $c
    $ws.Range("D$r").Value2 = $d
    $ws.Range("E$r").Value2 = $e
    $ws.Range("F$r").Value2 = $f
    $ws.Range("G$r").Value2 = $g
    $txtRng.Style = "Normal"

    $ws.Range("H$r").Value2 = $h
}

# ===========================================================================
# 1) Insert the new "2022-Q3" worksheet right after "总计" (position 2) and
#    populate it with the fund-holdings detail table.
# ===========================================================================

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Borrow the header / index-column formatting (bold, boxed, centered — style
# index "2" in the workbook) from the neighbouring quarter sheet so we don't
# introduce any brand-new cell styles. At this point in sheet order the
# "2022-Q2" sheet (formerly Sheet2) has been pushed to position 3.
$formatDonor = $wb.Worksheets.Item(3)
$formatDonor.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial($xlPasteFormats)
$formatDonor.Range("A2").Copy()
$q3Sheet.Range("A2:A23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Header row
$q3Sheet.Range("B1").Value2 = "基金代码"
$q3Sheet.Range("C1").Value2 = "基金名称"
$q3Sheet.Range("D1").Value2 = "基金规模"
$q3Sheet.Range("E1").Value2 = "股票总仓位"
$q3Sheet.Range("F1").Value2 = "仓位占比"
$q3Sheet.Range("G1").Value2 = "持有市值(亿元)"
$q3Sheet.Range("H1").Value2 = "仓位排名"

# Data rows
Set-Q3Row $q3Sheet 2  0  '004997' '广发高端制造股票A'                     '130.03' '87.00' '4.09' '5.3182' 9
Set-Q3Row $q3Sheet 3  1  '002132' '广发鑫享灵活配置混合A'                  '65.56'  '88.72' '3.41' '2.2356' 10
Set-Q3Row $q3Sheet 4  2  '014725' '广发成长动力三年持有期混合A'              '29.73'  '60.07' '1.72' '0.5114' 10
Set-Q3Row $q3Sheet 5  3  '010160' '广发高端制造股票C'                     '7.66'   '87.00' '4.09' '0.3133' 9
Set-Q3Row $q3Sheet 6  4  '015322' '广发鑫享灵活配置混合C'                  '5.01'   '88.72' '3.41' '0.1708' 10
Set-Q3Row $q3Sheet 7  5  '233006' '大摩领先优势混合'                      '3.47'   '94.24' '4.01' '0.1391' 10
Set-Q3Row $q3Sheet 8  6  '671030' '西部利得事件驱动股票'                   '2.90'   '94.40' '4.61' '0.1337' 10
Set-Q3Row $q3Sheet 9  7  '010322' '大摩新兴产业股票'                      '1.98'   '94.19' '4.70' '0.0931' 10
Set-Q3Row $q3Sheet 10 8  '014726' '广发成长动力三年持有期混合C'              '4.55'   '60.07' '1.72' '0.0783' 10
Set-Q3Row $q3Sheet 11 9  '010268' '太平睿安混合A'                        '4.03'   '39.63' '1.47' '0.0592' 7
Set-Q3Row $q3Sheet 12 10 '006973' '太平睿盈混合A'                        '3.84'   '28.79' '1.03' '0.0396' 2
Set-Q3Row $q3Sheet 13 11 '015043' '西部利得时代动力混合A'                   '0.68'   '79.48' '4.70' '0.0320' 5
Set-Q3Row $q3Sheet 14 12 '014571' '东吴安享量化灵活配置混合C'                '0.47'   '90.86' '6.02' '0.0283' 9
Set-Q3Row $q3Sheet 15 13 '580007' '东吴安享量化灵活配置混合A'                '0.47'   '90.86' '6.02' '0.0283' 9
Set-Q3Row $q3Sheet 16 14 '007669' '太平睿盈混合C'                        '1.04'   '28.79' '1.03' '0.0107' 2
Set-Q3Row $q3Sheet 17 15 '011886' '弘毅远方高端制造混合型发起式证券投资基金A'    '0.25'   '89.58' '3.55' '0.0089' 9
Set-Q3Row $q3Sheet 18 16 '015044' '西部利得时代动力混合C'                   '0.16'   '79.48' '4.70' '0.0075' 5
Set-Q3Row $q3Sheet 19 17 '001448' '华商双翼平衡混合'                      '0.49'   '39.66' '1.29' '0.0063' 9
Set-Q3Row $q3Sheet 20 18 '011887' '弘毅远方高端制造混合型发起式证券投资基金C'    '0.11'   '89.58' '3.55' '0.0039' 9
Set-Q3Row $q3Sheet 21 19 '015707' '安信新能源主题股票A'                    '0.18'   '53.26' '2.15' '0.0039' 9
Set-Q3Row $q3Sheet 22 20 '010269' '太平睿安混合C'                        '0.25'   '39.63' '1.47' '0.0037' 7
Set-Q3Row $q3Sheet 23 21 '015708' '安信新能源主题股票C'                    '0.14'   '53.26' '2.15' '0.0030' 9

# ===========================================================================
# 2) Update the "总计" summary sheet: a new 2022-Q3 row is inserted right
#    after the header, pushing every other quarter down by one row.
# ===========================================================================

$totalSheet.Range("A2").Value2 = 0
$totalSheet.Range("B2").Value2 = "2022-Q3"
$totalSheet.Range("C2").Value2 = 22
$totalSheet.Range("D2").Value2 = 9.23

$totalSheet.Range("A3").Value2 = 1
$totalSheet.Range("B3").Value2 = "2022-Q2"
$totalSheet.Range("C3").Value2 = 11
$totalSheet.Range("D3").Value2 = 9.81

$totalSheet.Range("A4").Value2 = 2
$totalSheet.Range("B4").Value2 = "2022-Q1"
$totalSheet.Range("C4").Value2 = 21
$totalSheet.Range("D4").Value2 = 8.22

$totalSheet.Range("A5").Value2 = 3
$totalSheet.Range("B5").Value2 = "2021-Q4"
$totalSheet.Range("C5").Value2 = 26
$totalSheet.Range("D5").Value2 = 15.31

$totalSheet.Range("A6").Value2 = 4
$totalSheet.Range("B6").Value2 = "2021-Q3"
$totalSheet.Range("C6").Value2 = 17
$totalSheet.Range("D6").Value2 = 2.75

$totalSheet.Range("A7").Value2 = 5
$totalSheet.Range("B7").Value2 = "2021-Q2"
$totalSheet.Range("C7").Value2 = 16
$totalSheet.Range("D7").Value2 = 5.05

# Row 8 is brand new — copy the index-column format (style "2") from A7
# before writing into it so it matches the rest of column A.
$totalSheet.Range("A7").Copy()
$totalSheet.Range("A8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$totalSheet.Range("A8").Value2 = 6
$totalSheet.Range("B8").Value2 = "2021-Q1"
$totalSheet.Range("C8").Value2 = 3
$totalSheet.Range("D8").Value2 = 0.15

Write-Output "2022-Q3 sheet added; 总计 summary updated."
